$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-28 15:05:32"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-28 15:05:27"
$wsZhCn.Range("K2").Value = "2016-08-28 15:05:51"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-28 15:05:32"
$wsDeDe.Range("K2").Value = "2016-08-28 15:05:57"
